# Update "Pais" worksheet with refreshed COVID-19 country data.
# Some countries swapped rank (and therefore rows) because totals changed;
# the country name column (A) and the numeric columns (B:H) are updated
# together per row to reflect the new ranking/values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $name, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Range("A$row").Value = $name
    $ws.Range("B$row").Value = $b
    $ws.Range("C$row").Value = $c
    $ws.Range("D$row").Value = $d
    $ws.Range("E$row").Value = $e
    $ws.Range("F$row").Value = $f
    $ws.Range("G$row").Value = $g
    $ws.Range("H$row").Value = $h
}

# Footer timestamp text (row 1)
$ws.Range("A1").Value = "Datos actualizados a 10 de Mayo de 2020 a las 12:04"

# Row 18 - Belgica (values updated, country unchanged)
Set-Row 18 "Belgica" 53081 485 13642 30783 476 75 8656

# Row 34 - Austria (values updated, country unchanged)
Set-Row 34 "Austria" 15871 38 13991 1262 72 3 618

# Row 37 - was Ucrania, now Rumania (new data)
Set-Row 37 "Rumania" 15362 231 7051 7359 242 13 952

# Row 38 - was Rumania, now Ucrania (unchanged data, shifted down a rank)
Set-Row 38 "Ucrania" 15232 522 3060 11781 201 15 391

# Row 39 - was Indonesia, now Banglades (new data)
Set-Row 39 "Banglades" 14657 887 2414 12015 1 14 228

# Row 40 - was Banglades, now Indonesia (unchanged data, shifted down a rank)
Set-Row 40 "Indonesia" 14032 387 2698 10361 0 14 973

# Row 55 - was Marruecos, now Finlandia (new data)
Set-Row 55 "Finlandia" 5913 33 4000 1648 45 0 265

# Row 56 - was Finlandia, now Marruecos (unchanged data, shifted down a rank)
Set-Row 56 "Marruecos" 5910 0 2461 3263 1 0 186

# Row 88 - was Eslovaquia, now Lituania (new data)
Set-Row 88 "Lituania" 1479 35 828 601 17 1 50

# Row 89 - was Eslovenia, now Eslovaquia (unchanged data, shifted up a rank)
Set-Row 89 "Eslovaquia" 1457 2 941 490 5 0 26

# Row 90 - was Lituania, now Eslovenia (unchanged data, shifted down a rank)
Set-Row 90 "Eslovenia" 1454 0 255 1098 10 0 101

# Row 93 - Hong Kong (values updated, country unchanged; F/G/H untouched)
$ws.Range("A93").Value = "Hong Kong"
$ws.Range("B93").Value = 1048
$ws.Range("C93").Value = 3
$ws.Range("D93").Value = 982
$ws.Range("E93").Value = 62

# Row 140 - was Cabo Verde, now Etiopia (new data)
Set-Row 140 "Etiopia" 239 29 99 135 1 0 5

# Row 141 - was Etiopia, now Cabo Verde (unchanged data, shifted down a rank)
Set-Row 141 "Cabo Verde" 236 0 56 178 0 0 2

# Row 192 - was Belice, now Nueva Caledonia (unchanged data, shifted up a rank; only D/H differ)
$ws.Range("A192").Value = "Nueva Caledonia"
$ws.Range("D192").Value = 18
$ws.Range("H192").Value = 0

# Row 193 - was Nueva Caledonia, now Belice (unchanged data, shifted down a rank; only D/H differ)
$ws.Range("A193").Value = "Belice"
$ws.Range("D193").Value = 16
$ws.Range("H193").Value = 2
